$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25, columns B,C,D,E,F,H,I,M,O
$data = @(
    @(12.54690179691499, 9.428555705286486, 3.966546793897766, 11.535794073821, 20.30730892550124, 7.344005520526261, 17.13428314545413, 14.54929968178754, 18.069893364678),
    @(11.86650553936505, 8.922758593140493, 3.910455146537271, 11.47162278415196, 20.31224683771509, 7.344005520526261, 17.27335154175093, 14.22167374363958, 18.14189883606025),
    @(11.42797464107968, 8.595601600096757, 3.875358465795483, 11.4373520887251, 20.32376456729335, 7.344005520526261, 17.36395814447174, 14.01934933536985, 18.1925679151638),
    @(11.24419623374694, 8.458185047238574, 3.860901676604904, 11.4246862178146, 20.33058498205351, 7.344005520526261, 17.4021913446437, 13.93672711513701, 18.21482996576221),
    @(11.21337845651545, 8.435122191506505, 3.858492108761521, 11.42266178314007, 20.33184571673729, 7.344005520526261, 17.40861902445337, 13.9230008245282, 18.21862374933678),
    @(11.42551646631549, 8.59376483089302, 3.875164108565561, 11.43717599934899, 20.32384794913827, 7.344005520526261, 17.36446846790231, 14.0182356058303, 18.19286162725208),
    @(12.31669087702581, 9.257641844380132, 3.947348735587823, 11.51260951617626, 20.30724673972444, 7.344005520526261, 17.18114978252266, 14.43664501186319, 18.09337485398065),
    @(13.89437944252291, 10.42524318532946, 4.083254698881213, 11.70068904893029, 20.34223137451643, 7.344005520526261, 16.86315111280813, 15.24294947752669, 17.94993617654534),
    @(14.94460345886195, 11.19862334985471, 4.17909443990309, 11.86242023282504, 20.40923989164066, 7.344005520526261, 16.6549334642763, 15.82037789357164, 17.87658233447673),
    @(15.3979629929915, 11.53175679950032, 4.221710528619367, 11.94085970524855, 20.44867126401433, 7.344005520526261, 16.56575881297943, 16.07859846550229, 17.85027171174593),
    @(15.56608437222396, 11.65519864166575, 4.237697966628852, 11.97123840520257, 20.46488404677578, 7.344005520526261, 16.53279100374435, 16.17564655711249, 17.84133124336472),
    @(15.53003523569678, 11.62873398901241, 4.234261613931228, 11.96466618142933, 20.46133548112117, 7.344005520526261, 16.53985554396108, 16.15477960181214, 17.84321112683374),
    @(15.41186594663051, 11.54196682452734, 4.22302888894238, 11.94334559077719, 20.44997946887235, 7.344005520526261, 16.56303046417862, 16.08659790346302, 17.84951563333662),
    @(15.33901949575846, 11.4884661357592, 4.216128676032954, 11.93037330717223, 20.44319018936531, 7.344005520526261, 16.57733015036179, 16.04473632765118, 17.8535107440574),
    @(14.91448010943389, 11.17647434534289, 4.176288827940787, 11.85738991654568, 20.406842587432, 7.344005520526261, 16.66087309395168, 15.80340475758274, 17.87844451041093),
    @(14.6477533168446, 10.98027504190382, 4.15159023928967, 11.81384674044101, 20.38683313349658, 7.344005520526261, 16.71354643489207, 15.6541458291898, 17.89555415004017),
    @(14.49204526988548, 10.86566894836808, 4.137292365557053, 11.78926157425524, 20.37616703969143, 7.344005520526261, 16.74436441596642, 15.56788195363555, 17.90605910556143),
    @(14.43893289096702, 10.82656422903293, 4.132435853305043, 11.78101713145296, 20.37270057656984, 7.344005520526261, 16.7548883520915, 15.53860633787325, 17.90972968774707),
    @(14.6763846301379, 11.00134286149211, 4.154229025037043, 11.81843457838212, 20.38887598053907, 7.344005520526261, 16.70788525317072, 15.67007825476693, 17.89366402545307),
    @(15.44667197253517, 11.56752610205544, 4.226332363879006, 11.94958984030701, 20.45328029826675, 7.344005520526261, 16.55620166963541, 16.1066451682839, 17.84763602965576),
    @(15.92935687388671, 11.92176558412576, 4.272575763138371, 12.03923133753449, 20.50283434154342, 7.344005520526261, 16.46173647047235, 16.38764722720541, 17.82351885165478),
    @(15.67365003035576, 11.73415234308457, 4.247978253242211, 11.9910375145202, 20.47570619317555, 7.344005520526261, 16.51172595601757, 16.23809566929803, 17.83584245701247),
    @(14.66344777067217, 10.99182373166203, 4.153036335554767, 11.81635901832843, 20.38794979985088, 7.344005520526261, 16.71044300595823, 15.66287660872981, 17.89451646932024),
    @(13.4863865831932, 10.12405879918569, 4.04715309641329, 11.6455957393317, 20.32551366280238, 7.344005520526261, 16.94472447365892, 15.02701620498741, 17.98315233507105)
)

$colIndices = @(2, 3, 4, 5, 6, 8, 9, 13, 15)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $colIndices.Count; $j++) {
        $ws.Cells.Item($rowNum, $colIndices[$j]).Value = $rowValues[$j]
    }
}
